$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-14 Wednesday" "2024-02-15 Thursday"

Replace-Text "806÷8=100, 6" "448÷6=74, 4"
Replace-Text "114÷5=22, 4" "952÷9=105, 7"
Replace-Text "898÷4=224, 2" "607÷6=101, 1"
Replace-Text "892÷2=446, 0" "536÷2=268, 0"
Replace-Text "929÷2=464, 1" "955÷6=159, 1"

Replace-Text "739÷4=184, 3" "643÷5=128, 3"
Replace-Text "939÷8=117, 3" "879÷6=146, 3"
Replace-Text "883÷8=110, 3" "609÷5=121, 4"
Replace-Text "346÷2=173, 0" "403÷8=50, 3"
Replace-Text "246÷5=49, 1" "144÷3=48, 0"

Replace-Text "925÷8=115, 5" "216÷7=30, 6"
Replace-Text "584÷2=292, 0" "908÷4=227, 0"
Replace-Text "641÷5=128, 1" "275÷2=137, 1"
Replace-Text "110÷8=13, 6" "833÷9=92, 5"
Replace-Text "924÷7=132, 0" "754÷6=125, 4"

Replace-Text "621÷7=88, 5" "187÷9=20, 7"
Replace-Text "196÷9=21, 7" "947÷6=157, 5"
Replace-Text "577÷9=64, 1" "162÷7=23, 1"
Replace-Text "496÷7=70, 6" "448÷4=112, 0"
Replace-Text "834÷4=208, 2" "838÷4=209, 2"

Replace-Text "878÷5=175, 3" "903÷9=100, 3"
Replace-Text "558÷3=186, 0" "163÷6=27, 1"
Replace-Text "490÷8=61, 2" "869÷4=217, 1"
Replace-Text "855÷6=142, 3" "708÷3=236, 0"
Replace-Text "956÷5=191, 1" "718÷4=179, 2"

Write-Output "Done"
